$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column C (y_0_forecast) and column E (y_1_forecast), rows 2-19
$values = @{
    2  = @{ C = 2.740959689118805;    E = 3.206168778303486 }
    3  = @{ C = -5.478010998490157;   E = -2.013762956649334 }
    4  = @{ C = -0.1091898317121864;  E = -1.305195642355683 }
    5  = @{ C = 3.371423250978856;    E = 0.8060632160631576 }
    6  = @{ C = 1.627570629117536;    E = 2.766358213445708 }
    7  = @{ C = -0.03183655677961861; E = 1.102200073559878 }
    8  = @{ C = 1.812248956008733;    E = 1.209672013646323 }
    9  = @{ C = 1.290465392296114;    E = 0.9879295308886871 }
    10 = @{ C = 1.57569012346459;     E = 1.643656926428561 }
    11 = @{ C = 1.73823635068906;     E = 1.765380623247137 }
    12 = @{ C = 2.337818484846466;    E = 2.076648015684435 }
    13 = @{ C = 0.8311911554373719;   E = 1.758956425699298 }
    14 = @{ C = -1.538034740964334;   E = -0.7351085756681308 }
    15 = @{ C = 0.303920243687994;    E = -0.6038293380915438 }
    16 = @{ C = 1.91914784107321;     E = 0.2384815980940092 }
    17 = @{ C = -0.1211988132392205;  E = 0.7707080878861294 }
    18 = @{ C = -0.04760886976447054; E = 0.5568966348730831 }
    19 = @{ C = 0.1245593350339691;   E = 0.1722027100061974 }
}

foreach ($row in $values.Keys) {
    $ws.Range("C$row").Value = $values[$row].C
    $ws.Range("E$row").Value = $values[$row].E
}
